$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we touch remain plain text (matches the source data,
# which stores prices/volumes/links as inline strings) so that Excel
# does not silently coerce strings like "0.9940" into the number 0.994.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.947.59'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.756.96'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9940'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.71%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.63'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -8.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9960'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.51%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5030'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -5.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.63'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -7.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2649'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -13.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06190'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -10.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.748.75'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.07%  '

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.68'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -14.81%  '

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.06927'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -11.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.495'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -10.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5947'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -21.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.77'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -13.47%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9916'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9945'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.988.48'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.71'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -16.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006785'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -14.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.969.23'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.083'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -11.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.095'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -13.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.125'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -14.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.96'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.533'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -9.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.862'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -15.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.90'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -12.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.69'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -7.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.778'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -11.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08085'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -8.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.474'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -14.94%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04485'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9934'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.622'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -10.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9995'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -11.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6077'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -16.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.699'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -13.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.960'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -15.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '104.91'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01527'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -11.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9946'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3851'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -19.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.164'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -12.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7363'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -19.09%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -10.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05212'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -10.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.994'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -20.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.27'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -13.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.56'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -12.92%  '
